$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each line: DateSerial,Hour,PredictionValue,LookupLabel
# Represents the post-edit content of worksheet rows 2..170 (one row per line,
# in order), after shifting the rolling 7-day forecast window forward so it
# now starts 18.02.2026 17:00 and ends 25.02.2026 17:00 (adds the
# "Necaluxan" and "Adrem" forecast days at the tail, drops the oldest two
# days from the head).
$data = @"
46071,17,0.288,18.02.202617
46071,18,0.12,18.02.202618
46071,19,0.015,18.02.202619
46071,20,0,18.02.202620
46071,21,0,18.02.202621
46071,22,0,18.02.202622
46071,23,0,18.02.202623
46071,24,0,18.02.202624
46072,1,0,19.02.20261
46072,2,0,19.02.20262
46072,3,0,19.02.20263
46072,4,0,19.02.20264
46072,5,0,19.02.20265
46072,6,0,19.02.20266
46072,7,0,19.02.20267
46072,8,0,19.02.20268
46072,9,0.215,19.02.20269
46072,10,1.188,19.02.202610
46072,11,1.808,19.02.202611
46072,12,1.822,19.02.202612
46072,13,1.963,19.02.202613
46072,14,1.933,19.02.202614
46072,15,1.79,19.02.202615
46072,16,1.535,19.02.202616
46072,17,0.5620000000000001,19.02.202617
46072,18,0.208,19.02.202618
46072,19,0.036,19.02.202619
46072,20,0,19.02.202620
46072,21,0,19.02.202621
46072,22,0,19.02.202622
46072,23,0,19.02.202623
46072,24,0,19.02.202624
46073,1,0,20.02.20261
46073,2,0,20.02.20262
46073,3,0,20.02.20263
46073,4,0,20.02.20264
46073,5,0,20.02.20265
46073,6,0,20.02.20266
46073,7,0,20.02.20267
46073,8,0,20.02.20268
46073,9,0.144,20.02.20269
46073,10,0.593,20.02.202610
46073,11,1.052,20.02.202611
46073,12,1.754,20.02.202612
46073,13,1.645,20.02.202613
46073,14,1.617,20.02.202614
46073,15,1.314,20.02.202615
46073,16,0.526,20.02.202616
46073,17,0.404,20.02.202617
46073,18,0.163,20.02.202618
46073,19,0.023,20.02.202619
46073,20,0,20.02.202620
46073,21,0,20.02.202621
46073,22,0,20.02.202622
46073,23,0,20.02.202623
46073,24,0,20.02.202624
46074,1,0,21.02.20261
46074,2,0,21.02.20262
46074,3,0,21.02.20263
46074,4,0,21.02.20264
46074,5,0,21.02.20265
46074,6,0,21.02.20266
46074,7,0,21.02.20267
46074,8,0,21.02.20268
46074,9,0.076,21.02.20269
46074,10,0.295,21.02.202610
46074,11,0.536,21.02.202611
46074,12,0.653,21.02.202612
46074,13,0.719,21.02.202613
46074,14,0.83,21.02.202614
46074,15,0.905,21.02.202615
46074,16,0.735,21.02.202616
46074,17,0.578,21.02.202617
46074,18,0.217,21.02.202618
46074,19,0.023,21.02.202619
46074,20,0,21.02.202620
46074,21,0,21.02.202621
46074,22,0,21.02.202622
46074,23,0,21.02.202623
46074,24,0,21.02.202624
46075,1,0,22.02.20261
46075,2,0,22.02.20262
46075,3,0,22.02.20263
46075,4,0,22.02.20264
46075,5,0,22.02.20265
46075,6,0,22.02.20266
46075,7,0,22.02.20267
46075,8,0,22.02.20268
46075,9,0.184,22.02.20269
46075,10,0.629,22.02.202610
46075,11,1.45,22.02.202611
46075,12,1.765,22.02.202612
46075,13,1.785,22.02.202613
46075,14,1.88,22.02.202614
46075,15,1.751,22.02.202615
46075,16,1.559,22.02.202616
46075,17,0.841,22.02.202617
46075,18,0.271,22.02.202618
46075,19,0.038,22.02.202619
46075,20,0,22.02.202620
46075,21,0,22.02.202621
46075,22,0,22.02.202622
46075,23,0,22.02.202623
46075,24,0,22.02.202624
46076,1,0,23.02.20261
46076,2,0,23.02.20262
46076,3,0,23.02.20263
46076,4,0,23.02.20264
46076,5,0,23.02.20265
46076,6,0,23.02.20266
46076,7,0,23.02.20267
46076,8,0,23.02.20268
46076,9,0.182,23.02.20269
46076,10,0.624,23.02.202610
46076,11,1.554,23.02.202611
46076,12,1.71,23.02.202612
46076,13,1.88,23.02.202613
46076,14,1.825,23.02.202614
46076,15,1.842,23.02.202615
46076,16,1.707,23.02.202616
46076,17,0.963,23.02.202617
46076,18,0.341,23.02.202618
46076,19,0.022,23.02.202619
46076,20,0,23.02.202620
46076,21,0,23.02.202621
46076,22,0,23.02.202622
46076,23,0,23.02.202623
46076,24,0,23.02.202624
46077,1,0,24.02.20261
46077,2,0,24.02.20262
46077,3,0,24.02.20263
46077,4,0,24.02.20264
46077,5,0,24.02.20265
46077,6,0,24.02.20266
46077,7,0,24.02.20267
46077,8,0,24.02.20268
46077,9,0.105,24.02.20269
46077,10,0.468,24.02.202610
46077,11,0.708,24.02.202611
46077,12,0.628,24.02.202612
46077,13,1.61,24.02.202613
46077,14,1.534,24.02.202614
46077,15,1.356,24.02.202615
46077,16,0.589,24.02.202616
46077,17,0.604,24.02.202617
46077,18,0.194,24.02.202618
46077,19,0.024,24.02.202619
46077,20,0,24.02.202620
46077,21,0,24.02.202621
46077,22,0,24.02.202622
46077,23,0,24.02.202623
46077,24,0,24.02.202624
46078,1,0,25.02.20261
46078,2,0,25.02.20262
46078,3,0,25.02.20263
46078,4,0,25.02.20264
46078,5,0,25.02.20265
46078,6,0,25.02.20266
46078,7,0,25.02.20267
46078,8,0,25.02.20268
46078,9,0.192,25.02.20269
46078,10,0.645,25.02.202610
46078,11,1.112,25.02.202611
46078,12,1.871,25.02.202612
46078,13,1.814,25.02.202613
46078,14,1.735,25.02.202614
46078,15,1.686,25.02.202615
46078,16,1.333,25.02.202616
46078,17,0.637,25.02.202617
"@

$lines = $data -split "`n" | Where-Object { $_.Trim() -ne "" }

$row = 2
foreach ($line in $lines) {
    $parts = $line.Trim() -split ","
    $dateSerial = [double]$parts[0]
    $hour = [double]$parts[1]
    $pred = [double]$parts[2]
    $lookup = [string]$parts[3]

    $ws.Cells.Item($row, 1).Value2 = $dateSerial
    $ws.Cells.Item($row, 2).Value2 = $hour
    $ws.Cells.Item($row, 3).Value2 = $pred
    $ws.Cells.Item($row, 4).Value = $lookup

    $row = $row + 1
}
